$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Users sheet: add "permissions" values (column F) and bump the
#    "modifiedAt" timestamps (column H) for the three existing users.
# ------------------------------------------------------------------
$usersWs = $wb.Worksheets.Item("Users")

$usersWs.Range("F2").Value = "read,write,delete,admin"
$usersWs.Range("H2").Value = "2025-05-01T00:34:45.831Z"

$usersWs.Range("F3").Value = "read,write,delete"
$usersWs.Range("H3").Value = "2025-05-01T00:34:53.271Z"

$usersWs.Range("F4").Value = "read"
$usersWs.Range("H4").Value = "2025-05-01T00:34:59.346Z"

# ------------------------------------------------------------------
# 2. AuditLog sheet: append six new audit rows (32-37) recording the
#    permission updates above.
# ------------------------------------------------------------------
$auditWs = $wb.Worksheets.Item("AuditLog")

$h32 = '{"before":{"id":"1","username":"admin","name":"Lisa Williams","email":"admin@example.com","role":"admin","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-29T04:58:50.090Z","status":"active"},"after":{"id":"1","username":"admin","name":"Lisa Williams","email":"admin@example.com","role":"admin","permissions":["read","write","delete","admin"],"modifiedBy":"admin","modifiedAt":"2025-05-01T00:34:45.831Z","status":"active"}}'
$h34 = '{"before":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-29T04:59:16.302Z","status":"active"},"after":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":["read","write","delete"],"modifiedBy":"admin","modifiedAt":"2025-05-01T00:34:53.271Z","status":"active"}}'
$h36 = '{"before":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-29T04:59:22.458Z","status":"active"},"after":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":["read"],"modifiedBy":"admin","modifiedAt":"2025-05-01T00:34:59.346Z","status":"active"}}'

$rows = @(
  @(32, "AUDIT1746059685832", "Users", "1", "UPDATE", "1", "admin", "2025-05-01T00:34:45.832Z", $h32, "Updated User 1"),
  @(33, "AUDIT1746059685832", "Users", "1", "UPDATE", "1", "admin", "2025-05-01T00:34:45.832Z", $h32, "Updated User 1"),
  @(34, "AUDIT1746059693271", "Users", "2", "UPDATE", "1", "admin", "2025-05-01T00:34:53.271Z", $h34, "Updated User 2"),
  @(35, "AUDIT1746059693271", "Users", "2", "UPDATE", "1", "admin", "2025-05-01T00:34:53.271Z", $h34, "Updated User 2"),
  @(36, "AUDIT1746059699346", "Users", "3", "UPDATE", "1", "admin", "2025-05-01T00:34:59.346Z", $h36, "Updated User 3"),
  @(37, "AUDIT1746059699346", "Users", "3", "UPDATE", "1", "admin", "2025-05-01T00:34:59.346Z", $h36, "Updated User 3")
)

# Columns C (entityId) and E (userId) hold digit-only strings ("1", "2",
# "3"...) for these particular rows - force text format on just those
# cells before writing so they don't get auto-coerced to numbers (every
# column on this sheet is text, per the source data).
$textCols = @(3, 5)

foreach ($row in $rows) {
    $r = $row[0]

    for ($col = 1; $col -le 9; $col++) {
        $cell = $auditWs.Cells.Item($r, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$col]
    }
}
